# RHI plots with the position of the sphere added
# Adds two new columns (AA: "Exp Constant", AB: "Exp Constant [dB]") to the
# "tabla" worksheet, filling the same constant value down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$expConstant = 385250961.9682089
$expConstantDb = 85.85743731821252

# Copy the header formatting (bold font + border + centered alignment) from
# the last existing header cell (Z1) onto the two new header cells so the
# new columns look consistent with the rest of the table.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (row 1)
$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"

# Data rows (row 2 through 40)
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 27).Value = $expConstant
    $ws.Cells.Item($r, 28).Value = $expConstantDb
}
